$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.783599999999992
$ws.Range("B8").Value = 4.991499999999998
$ws.Range("A12").Value = -22.5975
$ws.Range("B12").Value = 5.828100000000001
$ws.Range("B14").Value = 8.8012
$ws.Range("B22").Value = 4.808200000000004
